# Trello fetch/export feature: add a "Short ID" column (F) with the short
# Trello board ids, restyle the data table, resize the columns and move the
# active selection, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Short ID" column (F) -------------------------------------
# Written in this exact order so the shared-string table grows with the same
# new entries, in the same order, as the authored workbook:
#   14 "Short ID", 15 "3aSxVQ8E", 16 "wZlfoT6y", 17 "1DySJs3m"
$ws.Range("F1").Value = "Short ID"
$ws.Range("F3").Value = "3aSxVQ8E"
$ws.Range("F2").Value = "wZlfoT6y"
$ws.Range("F4").Value = "1DySJs3m"

# --- Apply the (no-op) alignment style used by the table ------------------
# This introduces the second cellXfs entry (applyAlignment="1") and assigns
# it to every cell of the table, A1:F4.
$ws.Range("A1:F4").WrapText = $false

# --- Resize the columns to match the authored layout -----------------------
$ws.Columns.Item(1).ColumnWidth = 15.142857
$ws.Columns.Item(2).ColumnWidth = 12.714286
$ws.Columns.Item(3).ColumnWidth = 9.285714
$ws.Columns.Item(4).ColumnWidth = 5.000000
$ws.Columns.Item(5).ColumnWidth = 18.000000
$ws.Columns.Item(6).ColumnWidth = 4.142857

# --- Move the active selection to H4, as left by the author ----------------
$null = $ws.Range("H4").Select()
